$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: AssembleWarship (sheet1) - insert a "MainType" column after "Type"
# ---------------------------------------------------------------------------
$wsShip = $wb.Worksheets.Item("AssembleWarship")
$wsShip.Columns("C").Insert()
$wsShip.Range("C1").Value = "MainType"
$wsShip.Range("C2").Value = "Warship"
$wsShip.Range("C3").Value = "Warship"
[void]$wsShip.Range("G14").Select()

# ---------------------------------------------------------------------------
# Sheet: AssembleWarShipType (sheet2) - header row gets the bold/fill style
# ---------------------------------------------------------------------------
$wsShipType = $wb.Worksheets.Item("AssembleWarShipType")
$wsShip.Range("A1:C1").Copy()
$wsShipType.Range("A1:C1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
[void]$wsShipType.Range("C24").Select()

# ---------------------------------------------------------------------------
# Sheet: AssembleWarshipClass (sheet3) - selection only change
# ---------------------------------------------------------------------------
$wsClass = $wb.Worksheets.Item("AssembleWarshipClass")
[void]$wsClass.Range("C24").Select()

# ---------------------------------------------------------------------------
# Sheet: AssembleParts (sheet4) - insert PartIconSmall/PartSprite columns
# after ModelTypeID, and an AssembleType column before CustomData
# ---------------------------------------------------------------------------
$wsParts = $wb.Worksheets.Item("AssembleParts")
$wsParts.Columns("C").Insert()
$wsParts.Columns("D").Insert()
$wsParts.Columns("G").Insert()

$wsParts.Range("C1").Value = "PartIconSmall"
$wsParts.Range("D1").Value = "PartSprite"
$wsParts.Range("G1").Value = "AssembleType"

$wsParts.Range("C2").Value = "SpriteOutput/Assemble/Icon/Assemble_Part_Engine_001"
$wsParts.Range("D2").Value = "SpriteOutput/Assemble/Icon/Assemble_Part_Engine_Icon_001"
$wsParts.Range("G2").Value = "WarShip"

$wsParts.Columns("C").ColumnWidth = 51.625
$wsParts.Columns("D").ColumnWidth = 51.625
$wsParts.Columns("G").ColumnWidth = 19.375

$wsParts.Activate()
[void]$wsParts.Range("G4").Select()

# ---------------------------------------------------------------------------
# Sheet: AssemblePartsType (sheet5) - untouched by this commit
# ---------------------------------------------------------------------------
